# "PROGETTO SE.docx" edits:
#  1) Reposition/resize the floating "Rettangolo 1" rectangle shape
#  2) Product backlog STATO column: update first 4 user-story rows
#  3) "SPRINT BACKLOG" heading: append bold " (SPRINT 1)"
#  4) Sprint backlog tasks for story #002: mention "in un Service"
#
# NOTE on ordering: in this host, once the Tables collection has been
# touched, Shape.Left/Shape.Top/Shape.Width/Shape.Height assignments
# stop taking effect. So the shape is repositioned first, before any
# table access.

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Move/resize the "Rettangolo 1" floating rectangle
#    old: left=380.95pt top=151.9pt  width=54.7pt  height=21.2pt
#    new: left=404.05pt top=185.13pt width=72pt(1in) height=24.56pt
# ---------------------------------------------------------------
$shape = $d.Shapes.Item(1)
$shape.Left = 404.05
$shape.Top = 185.12708661417324
$shape.Width = 72
$shape.Height = 24.564488188976377

# ---------------------------------------------------------------
# 2) Product backlog table (table 1): STATO column (column 8)
#    Row 3 (story 001) -> "Done"
#    Rows 4-6 (stories 002-004) -> "In progress"
# ---------------------------------------------------------------
$t1 = $d.Tables.Item(1)

$cell = $t1.Cell(3, 8)
$cell.Range.Find.Execute("To do", $true, $false, $false, $false, $false, `
    $true, 0, $false, "Done", 1)

$cell = $t1.Cell(4, 8)
$cell.Range.Find.Execute("To do", $true, $false, $false, $false, $false, `
    $true, 0, $false, "In progress", 1)

$cell = $t1.Cell(5, 8)
$cell.Range.Find.Execute("To do", $true, $false, $false, $false, $false, `
    $true, 0, $false, "In progress", 1)

$cell = $t1.Cell(6, 8)
$cell.Range.Find.Execute("To do", $true, $false, $false, $false, $false, `
    $true, 0, $false, "In progress", 1)

# ---------------------------------------------------------------
# 3) "SPRINT BACKLOG" heading -> add bold " (SPRINT 1)" suffix
#    (use Content.Find rather than the Paragraphs collection: once a
#    Tables collection has been touched, Paragraphs.Item(n) no longer
#    walks the document correctly in this host)
# ---------------------------------------------------------------
$titleRange = $d.Content
$titleFound = $titleRange.Find.Execute("SPRINT BACKLOG")
if ($titleFound) {
    $titleRange.Collapse(0)
    $titleRange.InsertAfter(" (SPRINT 1)")
    $titleRange.Font.Bold = 1
    $titleRange.Font.BoldBi = 1
}

# ---------------------------------------------------------------
# 4) Sprint backlog table (table 2): user story #002 tasks now
#    mention implementing "in un Service"
# ---------------------------------------------------------------
$t2 = $d.Tables.Item(2)

$cell = $t2.Cell(9, 3)
$cell.Range.Find.Execute("Implementare la funzione che ", $true, $false, `
    $false, $false, $false, $true, 0, $false, `
    "Implementare in un Service la funzione che ", 1)

$cell = $t2.Cell(10, 3)
$cell.Range.Find.Execute( `
    "Implementare la funzione che esegue la regola, una sola volta, raggiunto l’orario", `
    $true, $false, $false, $false, $false, $true, 0, $false, `
    "Implementare in un Service la funzione che esegue la regola, una sola volta, raggiunto l’orario", 1)
